$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Income side: add "NA" row at I6/J6
$ws.Range("I6").Value = "NA"
$ws.Range("J6").Formula = '=SUMIF($O:$O,I6,$N:$N)'

# Expenses side: add "NA" row at A17/B17 (copy formatting from the row above)
$ws.Range("A17").Value = "NA"
$ws.Range("B17").Formula = '=SUMIF($G:$G,A17,$F:$F)'

$ws.Range("A16").Copy() | Out-Null
$ws.Range("A17").PasteSpecial(-4122) | Out-Null

# Update the totals formula to include the new row
$ws.Range("A3").Formula = '=SUM(B5:B17)'

# Update selection to match diff
$ws.Range("J6").Select() | Out-Null
